$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 4 corresponds to RQ2 - update the requisite text and its evaluation columns
$ws.Range("B4").Value = "RQ2:El sistema debe permitir almacenar la información del cliente ya ingresada en la computadora"
$ws.Range("C4").Value = "si"
$ws.Range("D4").Value = "si"
$ws.Range("E4").Value = "si"
$ws.Range("F4").Value = "si"
$ws.Range("G4").Value = "no"
$ws.Range("H4").Value = "si "
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "si"

# Adjust the row height for row 4 to fit the new text
$ws.Rows.Item(4).RowHeight = 49.5

# J4 now evaluates to "si" (pass) - match the green fill used by the other
# passing rows (copy format from J3, which already has the "si" styling)
$ws.Range("J3").Copy()
$ws.Range("J4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the view: scroll back to the top and select B4
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B4").Select()
